$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.071.64'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.835.56'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.99'
$ws.Range('E5').Value = '  +1.75%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6364'
$ws.Range('E6').Value = '  +2.32%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07588'
$ws.Range('E8').Value = '  +2.82%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2953'
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.91'
$ws.Range('E10').Value = '  +1.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07756'
$ws.Range('E11').Value = '  +2.14%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.848.18'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.007'
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6713'
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '83.35'
$ws.Range('E15').Value = '  +1.56%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.000009722'
$ws.Range('E16').Value = '  +6.94%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.125'
$ws.Range('E17').Value = '  +1.87%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '29.087.23'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.60'
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '226.94'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  +0.62%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '160.44'
$ws.Range('E24').Value = '  +0.61%  '
$ws.Range('E25').Value = '  +3.48%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.550'
$ws.Range('E26').Value = '  +1.52%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.02'
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.502'
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.129'
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.077'
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.206'
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.05401'
$ws.Range('E32').Value = '  +2.94%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.864'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7492'
$ws.Range('E34').Value = '  +2.00%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.142'
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.662'
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.237.84'
$ws.Range('E37').Value = '  -3.56%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01794'
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.762'
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.646'
$ws.Range('E40').Value = '  +5.07%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9055'
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '102.30'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.981.19'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '65.07'
$ws.Range('E45').Value = '  +2.10%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00000000123'
$ws.Range('E46').Value = '  +2.32%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5110'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4091'
$ws.Range('E48').Value = '  +3.14%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.108'
$ws.Range('E49').Value = '  +3.41%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.773'
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05780'
$ws.Range('E51').Value = '  +0.53%  '
